$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add the new one ------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Power rails"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LS SE Gains"

# --- Fill in data on "LS SE Gains" -------------------------------------------
$ws2.Range("A3").Value = 38.15
$ws2.Range("B3").Value = 0.03947
$ws2.Range("C3").Value = 0.03964
$ws2.Range("D3").Formula = "=B3/A3*1000"
$ws2.Range("E3").Formula = "=C3/A3*1000"

$ws2.Range("A4").Value = 38.15
$ws2.Range("B4").Value = 0.38654
$ws2.Range("C4").Value = 0.38629
$ws2.Range("D4").Formula = "=B4/A4*1000"
$ws2.Range("E4").Formula = "=C4/A4*1000"

$ws2.Range("A5").Value = 38.15
$ws2.Range("B5").Value = 1.903
$ws2.Range("C5").Value = 1.89959
$ws2.Range("D5").Formula = "=B5/A5*1000"
$ws2.Range("E5").Formula = "=C5/A5*1000"

$ws2.Range("A6").Value = 38.15
$ws2.Range("B6").Value = 3.8672
$ws2.Range("C6").Value = 3.85968
$ws2.Range("D6").Formula = "=B6/A6*1000"
$ws2.Range("E6").Formula = "=C6/A6*1000"

# Header row
$ws2.Range("A2").Value = "Input (mV)"
$ws2.Range("B2").Value = "Hardware Output (V)"
$ws2.Range("C2").Value = "ADC averaged reading (V)"
$ws2.Range("D2").Value = "Hardware Gain (V/V)"
$ws2.Range("E2").Value = "Perceived Gain (V/V)"

# Methods column
$ws2.Range("G2").Value = "Methods"
$ws2.Range("G3").Value = "Input measured by multimeter in mV mode"
$ws2.Range("G4").Value = "Output measured by multimeter in mV or V mode"
$ws2.Range("G5").Value = "Output of ADC printed to serial terminal as well, averaged over 1000 samples measured at 50 Hz"

# Notes
$ws2.Range("A8").Value = "measured with OWON B35T+"
$ws2.Range("B8").Value = "measured with FLIR"
$ws2.Range("A9").Value = "through alligator cables"
$ws2.Range("B9").Value = "clipped with J hooks"

# --- Formatting ---------------------------------------------------------------
# data values (A3:C6) centered
$ws2.Range("A3:C6").HorizontalAlignment = -4108

# computed gains (D3:E6) centered with 5-decimal number format
$ws2.Range("D3:E6").NumberFormat = "0.00000"
$ws2.Range("D3:E6").HorizontalAlignment = -4108

# header row (A2:E2) bold, gray fill, centered
$ws2.Range("A2:E2").Font.Bold = $true
$ws2.Range("A2:E2").Interior.Color = 14277081
$ws2.Range("A2:E2").HorizontalAlignment = -4108

# "Methods" label keeps the existing bold-only style
$ws2.Range("G2").Font.Bold = $true

# --- Column widths -------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(3).ColumnWidth = 23.833333333333332
$ws2.Columns.Item(4).ColumnWidth = 18.0
$ws2.Columns.Item(5).ColumnWidth = 18.0

# --- Selections ------------------------------------------------------------
$ws1.Range("B26").Select()
$ws2.Range("F10").Select()
$ws2.Activate()
